$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "F27" = 63
    "G27" = 466.83
    "F36" = 8
    "G36" = 145.92
    "B47" = 15161.57
    "F80" = 143
    "G80" = 11607.31
    "F82" = 26
    "G82" = 1167.4
    "F83" = 25
    "G83" = 1669.75
    "F86" = 21
    "G86" = 2369.64
    "F103" = 14
    "G103" = 711.2
    "B114" = 267439.82
    "F193" = 304
    "G193" = 19699.2
    "F198" = 31
    "G198" = 2329.34
    "B200" = 48447.95
    "F219" = 5
    "G219" = 3360.2
    "B222" = 52259.56
    "B233" = 48719
    "E233" = 353.35
    "F233" = -81
    "G233" = -23955.75
    "B234" = 64979
    "E234" = 314.41
    "F234" = 11
    "G234" = 3253.25
    "B246" = 64973
    "E246" = 35.4
    "F246" = 64
    "G246" = 2131.2
    "B247" = 48706
    "E247" = 39.8
    "F247" = -144
    "G247" = -4795.2
    "F272" = 76
    "G272" = 6665.2
    "B274" = 90392.98
    "F276" = 21
    "G276" = 4540.62
    "B292" = 63520
    "E292" = 153.4
    "F292" = 73
    "G292" = 10532.44
    "B293" = 55373
    "E293" = 163.62
    "F293" = -94
    "G293" = -13562.32
    "B294" = 63571
    "E294" = 152.53
    "F294" = 4
    "G294" = 573.92
    "B295" = 57802
    "E295" = 162.71
    "F295" = -79
    "G295" = -11334.92
    "B299" = 55356
    "E299" = 54.04
    "F299" = -158
    "G299" = -7527.12
    "B300" = 63510
    "E300" = 50.66
    "F300" = 145
    "G300" = 6907.8
    "B311" = 61605
    "E311" = 133.78
    "F311" = -13
    "G311" = -1455.48
    "B312" = 63563
    "E312" = 119.04
    "F312" = 2
    "G312" = 223.92
    "B315" = 63560
    "E315" = 134.87
    "F315" = 1
    "G315" = 126.86
    "B316" = 60325
    "E316" = 151.57
    "F316" = -102
    "G316" = -12939.72
    "F321" = 105
    "G321" = 6163.5
    "F322" = 187
    "G322" = 19431.17
    "F333" = 719
    "G333" = 123186.27
    "B339" = 315499.67
    "F355" = 88
    "G355" = 2837.12
    "B356" = 31930
    "E356" = 26.8
    "F356" = -62
    "G356" = -1390.04
    "B357" = 63681
    "E357" = 23.84
    "F357" = 0
    "G357" = 0
    "B361" = 13003.44
    "F372" = 15
    "G372" = 809.25
    "B395" = 249837
    "B465" = 65069
    "E465" = 14.3
    "F465" = 2
    "G465" = 26.9
    "B466" = 53757
    "E466" = 16.08
    "F466" = -159
    "G466" = -2138.55
    "B469" = 53263
    "E469" = 15.29
    "F469" = -309
    "G469" = -3958.29
    "B470" = 65066
    "E470" = 13.61
    "F470" = 139
    "G470" = 1780.59
    "F476" = 117
    "G476" = 2308.41
    "B479" = 45718
    "E479" = 19.38
    "F479" = -294
    "G479" = -4768.68
    "B480" = 64927
    "E480" = 17.26
    "F480" = 198
    "G480" = 3211.56
    "F484" = 534
    "G484" = 3465.66
    "B485" = 45709
    "E485" = 15.69
    "F485" = -300
    "G485" = -3945
    "B486" = 64925
    "E486" = 13.97
    "F486" = 203
    "G486" = 2669.45
    "B490" = 65067
    "E490" = 15.65
    "F490" = 246
    "G490" = 3623.58
    "B491" = 53595
    "E491" = 17.61
    "F491" = -335
    "G491" = -4934.55
    "B492" = -2321.98
    "F495" = 19
    "G495" = 594.7
    "B508" = 12986.05
    "B576" = 64810
    "E576" = 291.22
    "F576" = 6
    "G576" = 1643.52
    "B577" = 53319
    "E577" = 310.64
    "F577" = -6
    "G577" = -1643.52
    "B608" = 64830
    "E608" = 34.9
    "F608" = 112
    "G608" = 3676.96
    "B609" = 60022
    "E609" = 37.22
    "F609" = -113
    "G609" = -3709.79
    "F662" = 348
    "G662" = 27972.24
    "B663" = 36996.64
    "F715" = 76
    "G715" = 7736.04
    "F719" = 93
    "G719" = 2019.96
    "F726" = 18
    "G726" = 673.5599999999999
    "F727" = 280
    "G727" = 33798.8
    "B728" = 147424.51
    "F736" = 28
    "G736" = 4211.48
    "F744" = 54
    "G744" = 2206.98
    "F747" = 108
    "G747" = 3570.48
    "F750" = 143
    "G750" = 7134.27
    "F753" = 122
    "G753" = 29499.6
    "B755" = 84703.5
    "F780" = 3201
    "G780" = 522115.11
    "F783" = 484
    "G783" = 70010.60000000001
    "B787" = 796464.5
    "F790" = 109
    "G790" = 15915.09
    "B804" = 79637.03
    "B805" = 3107817.4
    "B806" = 3107817.4
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
